$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted ahead of the existing row 115,
# pushing every subsequent record (old rows 115-167) down by one row
# (new rows 116-168), and the sheet's used range grew to A1:R168.
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new "Choclo" observation.
$ws.Range("A115").Value = 5
$ws.Range("B115").Value = "Macroferia Regional de Talca"
$ws.Range("C115").Value = "Maule"
$ws.Range("D115").Value = 44572
$ws.Range("E115").Value = 7
$ws.Range("F115").Value = 100112024
$ws.Range("G115").Value = "Choclo"
$ws.Range("H115").Value = "Choclero"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 50000
$ws.Range("K115").Value = 220
$ws.Range("L115").Value = 250
$ws.Range("M115").Value = 238
$ws.Range("N115").Value = "`$/unidad"
$ws.Range("O115").Value = "Región del Maule"
$ws.Range("P115").Value = 238
$ws.Range("Q115").Value = 1
$ws.Range("R115").Value = "Hortaliza"
